# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps to reflect the latest report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 625a38d0... row
$wsOverview.Range("G4").Value = "2016-08-31 06:51:25"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H4").Value = "2016-08-31 06:51:20"
$wsZhCn.Range("K4").Value = "2016-08-31 06:51:38"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H4").Value = "2016-08-31 06:51:25"
$wsDeDe.Range("K4").Value = "2016-08-31 06:51:46"
